$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.019.48"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "2.924.87"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'590.95"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").Value = "'146.80"
$ws.Range("E6").Value = "  +1.33%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").Value = "'6.90"
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("E10").Value = "  -0.23%  "
$ws.Range("D11").Value = "'0.442"
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("D13").Value = "'33.68"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("D15").Value = "3.409.66"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "60.977.33"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").Value = "2.925.69"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").Value = "'432.04"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("E20").Value = "  -1.54%  "
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").Value = "'81.40"
$ws.Range("E23").Value = "  +1.43%  "
$ws.Range("D24").Value = "'10.92"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "'11.92"
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  +5.04%  "
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("D30").Value = "'7.02"
$ws.Range("E30").Value = "  -2.33%  "
$ws.Range("D31").Value = "'26.68"
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("E32").Value = "  +2.56%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("D35").Value = "'1.02"
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").Value = "'3.01"
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("E38").Value = "  -0.80%  "
$ws.Range("E39").Value = "  -4.79%  "
$ws.Range("D40").Value = "'8.58"
$ws.Range("E40").Value = "  -0.79%  "
$ws.Range("B41").Value = "Arweave"
$ws.Range("C41").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D41").Value = "'41.27"
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").Value = "'0.283"
$ws.Range("E42").Value = "  -3.18%  "
$ws.Range("D43").Value = "'379.47"
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("D44").Value = "2.701.99"
$ws.Range("E44").Value = "  +1.13%  "
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("D46").Value = "'134.12"
$ws.Range("E46").Value = "  +1.42%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").Value = "'23.84"
$ws.Range("E48").Value = "  -2.01%  "
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("E50").Value = "  -2.28%  "
$ws.Range("E51").Value = "  -0.30%  "

# Reset style on forced-text cells so no stray style index is introduced
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
